# Insert a new weekly price record as row 37 on the single data sheet.
# Every existing row from 37 downward (through the old last row 113)
# shifts down by one row to make room; the sheet's used range grows
# from A1:T113 to A1:T114.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 37 (and everything below it) down by one row, carrying the
# existing formatting (e.g. the date style on column D) along with it.
$ws.Rows.Item(37).Insert()

# Populate the newly-opened row 37 with the new observation.
$ws.Cells.Item(37, 1).Value  = 7
$ws.Cells.Item(37, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(37, 3).Value  = "Ñuble"
$ws.Cells.Item(37, 4).Value  = 45014
$ws.Cells.Item(37, 5).Value  = 16
$ws.Cells.Item(37, 6).Value  = "Fruta"
$ws.Cells.Item(37, 7).Value  = 100108
$ws.Cells.Item(37, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(37, 9).Value  = 100108002
$ws.Cells.Item(37, 10).Value = "Mango"
$ws.Cells.Item(37, 11).Value = "Sin especificar"
$ws.Cells.Item(37, 12).Value = "Primera"
$ws.Cells.Item(37, 13).Value = 50
$ws.Cells.Item(37, 14).Value = 7500
$ws.Cells.Item(37, 15).Value = 8000
$ws.Cells.Item(37, 16).Value = 7800
$ws.Cells.Item(37, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(37, 18).Value = "Perú"
$ws.Cells.Item(37, 19).Value = 1950
$ws.Cells.Item(37, 20).Value = 4
